$wb = $excel.ActiveWorkbook

# Sheet "展览": F2, F12, F14, F16, F23, F26, F31 each +1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6565
$ws1.Range("F12").Value = 95
$ws1.Range("F14").Value = 497
$ws1.Range("F16").Value = 1037
$ws1.Range("F23").Value = 201
$ws1.Range("F26").Value = 145
$ws1.Range("F31").Value = 678

# Sheet "演出": F18 +1
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 325

# Sheet "全部类型": F8, F25, F29, F44, F51 each +1
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 6565
$ws4.Range("F25").Value = 497
$ws4.Range("F29").Value = 1037
$ws4.Range("F44").Value = 145
$ws4.Range("F51").Value = 678
